# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
# Swap the match-data columns (B:AB) between specific pairs of adjacent
# rows in the "Portugal Primeira Liga" sheet. Column A (row index) is
# left untouched; only the match record fields (id, HomeTeam, AwayTeam,
# scores, odds, etc.) are exchanged between the two rows of each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(76, 77),
    @(128, 129),
    @(133, 134),
    @(139, 140),
    @(151, 152),
    @(164, 165),
    @(258, 259),
    @(305, 306)
)

$firstCol = 2   # column B
$lastCol  = 28  # column AB

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range($ws.Cells.Item($r1, $firstCol), $ws.Cells.Item($r1, $lastCol))
    $range2 = $ws.Range($ws.Cells.Item($r2, $firstCol), $ws.Cells.Item($r2, $lastCol))

    # Capture current values of both rows before overwriting either.
    $values1 = @()
    $values2 = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $values1 += , ($ws.Cells.Item($r1, $c).Value2)
        $values2 += , ($ws.Cells.Item($r2, $c).Value2)
    }

    # Write swapped values back.
    $idx = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r1, $c).Value2 = $values2[$idx]
        $ws.Cells.Item($r2, $c).Value2 = $values1[$idx]
        $idx++
    }
}
